$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set E10:E15 (in_service) from FALSE to TRUE
$ws.Range("E10:E15").Value = $true
